# Split two run-on paragraphs in LOB1210.docx into multiple sentences
# separated by manual line breaks (<w:br/>), matching the target OOXML.
#
# Word's Find/Replace treats "^l" in replacement text as a manual line
# break (w:br), so a single Execute() both finds the run-on text and
# rewrites it with the breaks inserted at the right spots.

$d = $word.ActiveDocument

# --- "Critério" paragraph (Avaliação section) -----------------------------
$criterioOld = "O aluno poderá optar por dois critérios de avaliação:Critério 1: NF = (P1+P2)/2; ouCritério 2: NF = (NOTA 1 + NOTA 2)/2Sendo P1 e P2 avaliações escritas e NOTA 1 e NOTA 2 obtidas em atividades desenvolvidas em aula, trabalhos e relatórios de aulas práticas."
$criterioNew = "O aluno poderá optar por dois critérios de avaliação:^lCritério 1: NF = (P1+P2)/2; ou^lCritério 2: NF = (NOTA 1 + NOTA 2)/2^lSendo P1 e P2 avaliações escritas e NOTA 1 e NOTA 2 obtidas em atividades desenvolvidas em aula, trabalhos e relatórios de aulas práticas."

$d.Content.Find.Execute($criterioOld, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $criterioNew, 2)

# --- Bibliografia paragraph ------------------------------------------------
$biblioOld = "Bibliografia básica:1. LEPSCH, I.F. 19 Lições de pedologia. São Paulo, Oficina do Texto. 456p. 2011. ISBN 978-85-7975-029-8.Bibliografia complementar:1. CAMARGO, O.A. de; MONIZ, A.C.; JORGE, J.A.; VALADARES, J.M.A.S. Métodos de analise química, mineralógica e física de solos do Instituto Agronômico de Campinas. Campinas, Instituto Agronômico, 2009. 77 p. (Boletim técnico, 106, Edição revista e atualizada).2. DIAS Jr., M.S. Compactação do solo. In: Tópicos em ciência do solo, v.1. NOVAIS, R.F.; ALVAREZ, V.H.; SCHAEFER, G.R. (Eds.). Viçosa: SBCS, 2000. p.55-94.3. EMBRAPA – EMPRESA BRASILEIRA DE PESQUISA AGROPECUÁRIA. Manual de análises químicas de solos, plantas e fertilizantes. SILVA, F. C. da (org.). EMBRAPA Comunicação para Transferência de Tecnologia. Brasília: EMBRAPA, Solos, 1999b. 370p.4. OLIVEIRA, J.B. Pedologia aplicada. 3a Edição. Piracicaba: Ed. FEALQ, 2008. 592p.5. REICHARDT, K.; TIMM, L.C. Solo, planta e atmosfera: conceitos, processos e aplicações. Barueri: SP: ed. Manole, 2004. 478p.6. SCHNEIDER, P.; GIASSON, E.; KLAMT, E. Classificação da aptidão agrícola das terras: um sistema alternativo. Porto Alegre: UFRGS, 2007. 72p."
$biblioNew = "Bibliografia básica:^l1. LEPSCH, I.F. 19 Lições de pedologia. São Paulo, Oficina do Texto. 456p. 2011. ISBN 978-85-7975-029-8.^l^lBibliografia complementar:^l1. CAMARGO, O.A. de; MONIZ, A.C.; JORGE, J.A.; VALADARES, J.M.A.S. Métodos de analise química, mineralógica e física de solos do Instituto Agronômico de Campinas. Campinas, Instituto Agronômico, 2009. 77 p. (Boletim técnico, 106, Edição revista e atualizada).^l2. DIAS Jr., M.S. Compactação do solo. In: Tópicos em ciência do solo, v.1. NOVAIS, R.F.; ALVAREZ, V.H.; SCHAEFER, G.R. (Eds.). Viçosa: SBCS, 2000. p.55-94.^l3. EMBRAPA – EMPRESA BRASILEIRA DE PESQUISA AGROPECUÁRIA. Manual de análises químicas de solos, plantas e fertilizantes. SILVA, F. C. da (org.). EMBRAPA Comunicação para Transferência de Tecnologia. Brasília: EMBRAPA, Solos, 1999b. 370p.^l4. OLIVEIRA, J.B. Pedologia aplicada. 3a Edição. Piracicaba: Ed. FEALQ, 2008. 592p.^l5. REICHARDT, K.; TIMM, L.C. Solo, planta e atmosfera: conceitos, processos e aplicações. Barueri: SP: ed. Manole, 2004. 478p.^l6. SCHNEIDER, P.; GIASSON, E.; KLAMT, E. Classificação da aptidão agrícola das terras: um sistema alternativo. Porto Alegre: UFRGS, 2007. 72p."

$d.Content.Find.Execute($biblioOld, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $biblioNew, 2)
